$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Volume number text run: "49" -> "50"
$ws.Range("A8").Characters(21, 2).Text = "50"

# Week-covering date text runs: "12/5/2022" -> "12/12/2022", "12/11/2022" -> "12/18/2022"
# (apply right-to-left so earlier character offsets stay valid after length changes)
$ws.Range("C9").Characters(47, 10).Text = "12/18/2022"
$ws.Range("C9").Characters(27, 9).Text = "12/12/2022"

# Row 14
$ws.Range("C14").Value = 8
$ws.Range("D14").Value = 16
$ws.Range("E14").Value = -50
$ws.Range("F14").Value = 31
$ws.Range("G14").Value = 41
$ws.Range("H14").Value = -24.390243902439
$ws.Range("I14").Value = 414
$ws.Range("J14").Value = 473
$ws.Range("K14").Value = -12.473572938689
$ws.Range("L14").Value = -7.589285714285
$ws.Range("M14").Value = -20.231213872832
$ws.Range("N14").Value = -77.741935483871

# Row 15
$ws.Range("C15").Value = 34
$ws.Range("D15").Value = 33
$ws.Range("E15").Value = 3.030303030303
$ws.Range("F15").Value = 96
$ws.Range("G15").Value = 121
$ws.Range("H15").Value = -20.661157024793
$ws.Range("I15").Value = 1567
$ws.Range("J15").Value = 1456
$ws.Range("K15").Value = 7.623626373626
$ws.Range("L15").Value = 11.134751773049
$ws.Range("M15").Value = 17.290419161676
$ws.Range("N15").Value = -50.111429481057

# Row 16
$ws.Range("C16").Value = 284
$ws.Range("D16").Value = 312
$ws.Range("E16").Value = -8.974358974358
$ws.Range("F16").Value = 1202
$ws.Range("G16").Value = 1279
$ws.Range("H16").Value = -6.020328381548
$ws.Range("I16").Value = 16868
$ws.Range("J16").Value = 13313
$ws.Range("K16").Value = 26.703222414181
$ws.Range("L16").Value = 32.996925017740
$ws.Range("M16").Value = -10.666242982734
$ws.Range("N16").Value = -79.668779982161

# Row 17
$ws.Range("C17").Value = 397
$ws.Range("D17").Value = 444
$ws.Range("E17").Value = -10.585585585585
$ws.Range("F17").Value = 1769
$ws.Range("G17").Value = 1772
$ws.Range("H17").Value = -0.169300225733
$ws.Range("I17").Value = 25123
$ws.Range("J17").Value = 22340
$ws.Range("K17").Value = 12.457475380483
$ws.Range("L17").Value = 24.983831650166
$ws.Range("M17").Value = 51.773092490787
$ws.Range("N17").Value = -37.238002448225

# Row 18
$ws.Range("C18").Value = 272
$ws.Range("D18").Value = 285
$ws.Range("E18").Value = -4.561403508771
$ws.Range("F18").Value = 1164
$ws.Range("G18").Value = 1230
$ws.Range("H18").Value = -5.365853658536
$ws.Range("I18").Value = 15198
$ws.Range("J18").Value = 12266
$ws.Range("K18").Value = 23.903473014837
$ws.Range("L18").Value = 1.158146964856
$ws.Range("M18").Value = -16.130456376579
$ws.Range("N18").Value = -84.419019499292

# Row 19
$ws.Range("C19").Value = 993
$ws.Range("D19").Value = 1336
$ws.Range("E19").Value = -25.673652694610
$ws.Range("F19").Value = 3846
$ws.Range("G19").Value = 5155
$ws.Range("H19").Value = -25.392822502424
$ws.Range("I19").Value = 49853
$ws.Range("J19").Value = 39083
$ws.Range("K19").Value = 27.556738223780
$ws.Range("L19").Value = 44.154642454385
$ws.Range("M19").Value = 35.264271760364
$ws.Range("N19").Value = -39.880370946540

# Row 20
$ws.Range("C20").Value = 278
$ws.Range("D20").Value = 239
$ws.Range("E20").Value = 16.317991631799
$ws.Range("F20").Value = 1134
$ws.Range("G20").Value = 919
$ws.Range("H20").Value = 23.394994559303
$ws.Range("I20").Value = 13182
$ws.Range("J20").Value = 9996
$ws.Range("K20").Value = 31.872749099639
$ws.Range("L20").Value = 50.376454483230
$ws.Range("M20").Value = 31.714628297362
$ws.Range("N20").Value = -87.807427276511

# Row 21
$ws.Range("C21").Value = 2266
$ws.Range("D21").Value = 2665
$ws.Range("E21").Value = -14.971857410881
$ws.Range("F21").Value = 9242
$ws.Range("G21").Value = 10517
$ws.Range("H21").Value = -12.123229057716
$ws.Range("I21").Value = 122205
$ws.Range("J21").Value = 98927
$ws.Range("K21").Value = 23.530482072639
$ws.Range("L21").Value = 31.382035155620
$ws.Range("M21").Value = 19.486678073820
$ws.Range("N21").Value = -70.664416577047

# Row 22
$ws.Range("C22").Value = 44
$ws.Range("D22").Value = 52
$ws.Range("E22").Value = -15.384615384615
$ws.Range("F22").Value = 176
$ws.Range("G22").Value = 214
$ws.Range("H22").Value = -17.757009345794
$ws.Range("I22").Value = 2237
$ws.Range("J22").Value = 1736
$ws.Range("K22").Value = 28.859447004608
$ws.Range("L22").Value = 31.125439624853
$ws.Range("M22").Value = 6.018957345971

# Row 23
$ws.Range("C23").Value = 99
$ws.Range("D23").Value = 107
$ws.Range("E23").Value = -7.476635514018
$ws.Range("F23").Value = 414
$ws.Range("G23").Value = 456
$ws.Range("H23").Value = -9.210526315789
$ws.Range("I23").Value = 5737
$ws.Range("J23").Value = 5362
$ws.Range("K23").Value = 6.993659082431
$ws.Range("L23").Value = 14.717056588682
$ws.Range("M23").Value = 40.543851053405

# Row 24
$ws.Range("C24").Value = 2132
$ws.Range("D24").Value = 2017
$ws.Range("E24").Value = 5.701536936043
$ws.Range("F24").Value = 8861
$ws.Range("G24").Value = 8152
$ws.Range("H24").Value = 8.697252208047
$ws.Range("I24").Value = 111754
$ws.Range("J24").Value = 83435
$ws.Range("K24").Value = 33.941391502367
$ws.Range("L24").Value = 41.237282780410
$ws.Range("M24").Value = 41.333738032907

# Row 25
$ws.Range("C25").Value = 643
$ws.Range("D25").Value = 762
$ws.Range("E25").Value = -15.616797900262
$ws.Range("F25").Value = 2832
$ws.Range("G25").Value = 2991
$ws.Range("H25").Value = -5.315947843530
$ws.Range("I25").Value = 39828
$ws.Range("J25").Value = 35169
$ws.Range("K25").Value = 13.247462253689
$ws.Range("L25").Value = 23.885657407695
$ws.Range("M25").Value = -10.198191698045

# Row 26
$ws.Range("C26").Value = 49
$ws.Range("D26").Value = 52
$ws.Range("E26").Value = -5.769230769230
$ws.Range("F26").Value = 164
$ws.Range("G26").Value = 198
$ws.Range("H26").Value = -17.171717171717
$ws.Range("I26").Value = 2513
$ws.Range("J26").Value = 2392
$ws.Range("K26").Value = 5.058528428093
$ws.Range("L26").Value = 14.331210191082

# Row 27
$ws.Range("C27").Value = 79
$ws.Range("D27").Value = 95
$ws.Range("E27").Value = -16.842105263157
$ws.Range("F27").Value = 337
$ws.Range("G27").Value = 409
$ws.Range("H27").Value = -17.603911980440
$ws.Range("I27").Value = 5005
$ws.Range("J27").Value = 4794
$ws.Range("K27").Value = 4.401335002085
$ws.Range("L27").Value = 35.710412147505

# Row 28
$ws.Range("C28").Value = 22
$ws.Range("D28").Value = 38
$ws.Range("E28").Value = -42.105263157894
$ws.Range("F28").Value = 85
$ws.Range("G28").Value = 137
$ws.Range("H28").Value = -37.956204379562
$ws.Range("I28").Value = 1524
$ws.Range("J28").Value = 1825
$ws.Range("K28").Value = -16.493150684931
$ws.Range("L28").Value = -16.171617161716
$ws.Range("M28").Value = -11.395348837209
$ws.Range("N28").Value = -73.323997899527

# Row 29
$ws.Range("C29").Value = 17
$ws.Range("D29").Value = 33
$ws.Range("E29").Value = -48.484848484848
$ws.Range("F29").Value = 73
$ws.Range("G29").Value = 120
$ws.Range("H29").Value = -39.166666666666
$ws.Range("I29").Value = 1257
$ws.Range("J29").Value = 1523
$ws.Range("K29").Value = -17.465528562048
$ws.Range("L29").Value = -15.46738399462
$ws.Range("M29").Value = -11.789473684210
$ws.Range("N29").Value = -75.511396843950

# Row 30
$ws.Range("C30").Value = 2
$ws.Range("D30").Value = 12
$ws.Range("E30").Value = -83.333333333333
$ws.Range("F30").Value = 20
$ws.Range("G30").Value = 36
$ws.Range("H30").Value = -44.444444444444
$ws.Range("I30").Value = 602
$ws.Range("J30").Value = 518
$ws.Range("K30").Value = 16.216216216216
$ws.Range("L30").Value = 135.15625
